$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '39.974.37'
$ws.Range("E2").Value = '  +0.46%  '

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.236.23'
$ws.Range("E3").Value = '  -3.73%  '

$ws.Range("E4").Value = '  +0.00%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '293.32'
$ws.Range("E5").Value = '  -4.76%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '85.96'
$ws.Range("E6").Value = '  +3.77%  '

$ws.Range("E7").Value = '  -1.50%  '

$ws.Range("E8").Value = '  -0.02%  '

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.469'
$ws.Range("E9").Value = '  -1.12%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.0795'
$ws.Range("E10").Value = '  -0.21%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '30.40'
$ws.Range("E11").Value = '  +3.51%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '47.22'
$ws.Range("E12").Value = '  -9.82%  '

$ws.Range("E13").Value = '  -2.03%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.37'
$ws.Range("E14").Value = '  +0.91%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '2.585.15'
$ws.Range("E15").Value = '  -3.65%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '14.16'
$ws.Range("E16").Value = '  -2.65%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '2.227.44'
$ws.Range("E17").Value = '  -4.87%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.724'
$ws.Range("E18").Value = '  -2.99%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '39.892.05'
$ws.Range("E19").Value = '  +0.42%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.0₃0891'
$ws.Range("E20").Value = '  +0.30%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '5.79'
$ws.Range("E21").Value = '  -3.20%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '10.73'
$ws.Range("E22").Value = '  +3.60%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '65.42'
$ws.Range("E23").Value = '  -3.29%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '235.02'
$ws.Range("E24").Value = '  +0.88%  '

$ws.Range("E25").Value = '  -0.06%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.42'
$ws.Range("E26").Value = '  -3.19%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.84'
$ws.Range("E27").Value = '  +2.90%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '22.96'
$ws.Range("E28").Value = '  -0.95%  '

$ws.Range("E29").Value = '  +0.77%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.23'
$ws.Range("E30").Value = '  +1.14%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '33.36'
$ws.Range("E31").Value = '  -0.27%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '155.28'
$ws.Range("E32").Value = '  +1.76%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '4.85'
$ws.Range("E34").Value = '  -3.13%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0710'
$ws.Range("E35").Value = '  +0.47%  '

$ws.Range("E36").Value = '  -4.09%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '16.46'
$ws.Range("E37").Value = '  +8.45%  '

$ws.Range("E38").Value = '  -0.66%  '

$ws.Range("E39").Value = '  +2.26%  '

$ws.Range("E40").Value = '  -1.10%  '

$ws.Range("E41").Value = '  -0.35%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '3.76'
$ws.Range("E42").Value = '  +1.78%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.952.15'

$ws.Range("E44").Value = '  -3.07%  '

$ws.Range("E45").Value = '  +3.90%  '

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '9.50'
$ws.Range("E46").Value = '  +1.20%  '

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '16.26'
$ws.Range("E47").Value = '  -4.98%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.60'
$ws.Range("E48").Value = '  -0.85%  '

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '2.455.30'
$ws.Range("E49").Value = '  -3.86%  '

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '70.78'
$ws.Range("E50").Value = '  +2.28%  '

$ws.Range("E51").Value = '  +9.24%  '
